$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TrialData")
$design = $wb.Worksheets.Item("Design")

# ---------------------------------------------------------------------------
# Header row (row 1) - columns renamed/reordered
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "user"
$ws.Range("C1").Value = "user_text"
$ws.Range("D1").Value = "year"
$ws.Range("E1").Value = "field"
$ws.Range("F1").Value = "field_external_id"
$ws.Range("G1").Value = "inputs"
$ws.Range("H1").Value = "inputs_text"
$ws.Range("I1").Value = "latitude"
$ws.Range("J1").Value = "elevation"
$ws.Range("K1").Value = "soil_b"
$ws.Range("L1").Value = "weather_source"
$ws.Range("M1").Value = "weather_source_text"
$ws.Range("N1").Value = "weather_source_info"
$ws.Range("O1").Value = "date_sow"
$ws.Range("P1").Value = "date_harvest"
$ws.Range("Q1").Value = "date_emerg"
$ws.Range("R1").Value = "emergence"
$ws.Range("S1").Value = "plant_pop"
$ws.Range("T1").Value = "canopy1"
$ws.Range("U1").Value = "canopy2"
$ws.Range("V1").Value = "pop1"
$ws.Range("W1").Value = "pop2"
$ws.Range("X1").Value = "pop3"
$ws.Range("Y1").Value = "can1adoy"
$ws.Range("Z1").Value = "can1acov"
$ws.Range("AA1").Value = "can1atem"
$ws.Range("AB1").Value = "can1bdoy"
$ws.Range("AC1").Value = "can1bcov"
$ws.Range("AD1").Value = "can1btem"
$ws.Range("AE1").Value = "can1cdoy"
$ws.Range("AF1").Value = "can1ccov"
$ws.Range("AG1").Value = "can1ctem"
$ws.Range("AH1").Value = "can1ddoy"
$ws.Range("AI1").Value = "can1dcov"
$ws.Range("AJ1").Value = "can1dtem"
$ws.Range("AK1").Value = "can1edoy"
$ws.Range("AL1").Value = "can1ecov"
$ws.Range("AM1").Value = "can1etem"
$ws.Range("AN1").Value = "can2adoy"
$ws.Range("AO1").Value = "can2alos"
$ws.Range("AP1").Value = "can2bdoy"
$ws.Range("AQ1").Value = "can2blos"

# ---------------------------------------------------------------------------
# Data rows 2-10
# Column layout: A id(formula) B user C user_text(formula) D year E field
#                F field_external_id G inputs H inputs_text(formula) I latitude
#                J elevation K soil_b L weather_source M weather_source_text(formula)
#                N weather_source_info O date_sow P date_harvest Q date_emerg
#                R..AQ the remaining measurement columns (unchanged values)
# ---------------------------------------------------------------------------

$years    = @{2=2015;3=2016;4=2017;5=2018;6=2019;7=2020;8=2021;9=2022;10=2023}
$fields   = @{2="Ädelholm_2";3="Ädelholm_3";4="Ädelholm_4";5="Ädelholm_1";6="Ädelholm_2";7="Ädelholm_3";8="Ädelholm_4";9="Ädelholm_1";10="Ädelholm_2"}
$lats     = @{2="55.662999999999997";3="55.664999999999999";4="55.667000000000002";5="55.661000000000001";6="55.662999999999997";7="55.664999999999999";8="55.667000000000002";9="55.661000000000001";10="55.662999999999997"}
$sow      = @{2=42095;3=42461;4=42826;5=43191;6=43556;7=43922;8=44287;9=44645;10=45023}
$harvest  = @{2=42309;3=42675;4=43040;5=43405;6=43770;7=44136;8=44501;9=44866;10=45231}
$emerg    = @{2=42118;3=42485;4=42850;5=43215;6=43580;7=43946;8=44311;9=44659;10=45033}

for ($r = 2; $r -le 10; $r++) {
    $ws.Range("B$r").Value = 1
    $ws.Range("D$r").Value = $years[$r]
    $ws.Range("E$r").Value = $fields[$r]
    $ws.Range("F$r").Value = 282131
    $ws.Range("G$r").Value = 1
    $ws.Range("I$r").Value = $lats[$r]
    $ws.Range("J$r").Value = 20
    $ws.Range("K$r").Value = 2.4
    $ws.Range("L$r").Value = 1
    $ws.Range("N$r").Value = 40141
    $ws.Range("O$r").Value = $sow[$r]
    $ws.Range("P$r").Value = $harvest[$r]
    $ws.Range("Q$r").Value = $emerg[$r]

    $ws.Range("C$r").Formula = "=VLOOKUP(B$r,Design!F`$1:G`$1,2)"
    $ws.Range("H$r").Formula = "=VLOOKUP(G$r,Design!A`$1:B`$5,2)"
    $ws.Range("M$r").Formula = "=VLOOKUP(L$r,Design!C`$1:D`$2,2)"

    $ws.Range("R$r").Value = 0
    $ws.Range("S$r").Value = 0
    $ws.Range("T$r").Value = 0
    $ws.Range("U$r").Value = 0
    $ws.Range("V$r").Value = 90000
    $ws.Range("W$r").Value = 90000
    $ws.Range("X$r").Value = 90000
    $ws.Range("Y$r").Value = 0
    $ws.Range("Z$r").Value = 0
    $ws.Range("AA$r").Value = 0
    $ws.Range("AB$r").Value = 0
    $ws.Range("AC$r").Value = 0
    $ws.Range("AD$r").Value = 0
    $ws.Range("AE$r").Value = 0
    $ws.Range("AF$r").Value = 0
    $ws.Range("AG$r").Value = 0
    $ws.Range("AH$r").Value = 0
    $ws.Range("AI$r").Value = 0
    $ws.Range("AJ$r").Value = 0
    $ws.Range("AK$r").Value = 0
    $ws.Range("AL$r").Value = 0
    $ws.Range("AM$r").Value = 0
    $ws.Range("AN$r").Value = 0
    $ws.Range("AO$r").Value = 0
    $ws.Range("AP$r").Value = 0
    $ws.Range("AQ$r").Value = 0
}

# Column A ("id") - shared formula groups matching the original fill-down pattern:
#   A2:A7 one shared group, A8:A9 another, A10 its own standalone formula.
$ws.Range("A2:A7").Formula = "=B2*10000+G2+D2*100-200000"
$ws.Range("A8:A9").Formula = "=B8*10000+G8+D8*100-200000"
$ws.Range("A10").Formula = "=B10*10000+G10+D10*100-200000"

# ---------------------------------------------------------------------------
# Design sheet tweaks
# ---------------------------------------------------------------------------
$design.Range("G1").Value = "NBR"

# ---------------------------------------------------------------------------
# Column width / view cosmetics
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 10

# Update selection to match the final saved state
$ws.Range("A8").Select()

$excel.Calculation = -4105
$wb.Save()
